$d = $word.ActiveDocument

# 1) Resolve all tracked changes (anderson's g/G, c/C, "Datos del recurso"
#    rewordings, and Josue's deleted space) by accepting every revision in
#    the document. This mirrors the editor reviewing + accepting the
#    suggested changes.
$d.AcceptAllRevisions()

# 2) The accept-all above merges the old "algebraica,polinomio, monomio"
#    proofing run incorrectly (keeps the comma glued to "polinomio" and
#    leaves the spell-check bracket at the end). Patch the underlying OOXML
#    directly so the comma/space become their own runs and the spellcheck
#    proofErr bracket wraps only "algebraica,polinomio".
$flat = $d.WordOpenXML

$oldProof = '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>algebraica,polinomio,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>monomio</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$newProof = '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>algebraica,polinomio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>monomio</w:t></w:r>'
$flat = $flat.Replace($oldProof, $newProof)

# 3) Move the "_GoBack" bookmark: it used to sit around the "g" in the
#    title; Word has since relocated it (as Word does with the last-edit
#    position) to span the "Tratamiento de la informacion y competencia
#    digital" table row instead.
$oldBookmarkStart = '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$flat = $flat.Replace($oldBookmarkStart, '')
$oldBookmarkEnd = '<w:bookmarkEnd w:id="0"/>'
$flat = $flat.Replace($oldBookmarkEnd, '')

$cellAnchor = '<w:tcW w:w="4111" w:type="dxa"/></w:tcPr><w:p w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidRDefault="00044F81" w:rsidP="00CD26C8"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Tratamiento de la información y competencia digital</w:t></w:r>'
$newCell = '<w:tcW w:w="4111" w:type="dxa"/></w:tcPr><w:p w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidRDefault="00044F81" w:rsidP="00CD26C8"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack" w:colFirst="0" w:colLast="4"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Tratamiento de la información y competencia digital</w:t></w:r>'
$flat = $flat.Replace($cellAnchor, $newCell)

$rowAnchor = '<w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:tr w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidTr="00CD26C8"><w:tc><w:tcPr><w:tcW w:w="4536" w:type="dxa"/></w:tcPr><w:p w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidRDefault="00044F81" w:rsidP="00CD26C8"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> en comunicación'
$newRow = '<w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr></w:p></w:tc></w:tr><w:bookmarkEnd w:id="0"/><w:tr w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidTr="00CD26C8"><w:tc><w:tcPr><w:tcW w:w="4536" w:type="dxa"/></w:tcPr><w:p w:rsidR="00044F81" w:rsidRPr="00AC7496" w:rsidRDefault="00044F81" w:rsidP="00CD26C8"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>…</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> en comunicación'
$flat = $flat.Replace($rowAnchor, $newRow)

# 4) Drop the word/people.xml part (author presence-info list) along with
#    its content-type declaration and relationship, since the revisions
#    that referenced those authors were just accepted above.
$peoplePartPattern = '<pkg:part pkg:name="/word/people.xml"[^>]*><pkg:xmlData>.*?</pkg:xmlData></pkg:part>'
$flat = [System.Text.RegularExpressions.Regex]::Replace($flat, $peoplePartPattern, '')

$flat = $flat.Replace('<Override PartName="/word/people.xml" ContentType="application/vnd.openxmlformats-officedocument.wordprocessingml.people+xml"/>', '')
$flat = $flat.Replace('<Relationship Id="rId6" Type="http://schemas.microsoft.com/office/2011/relationships/people" Target="people.xml"/>', '')

$d.WordOpenXML = $flat

Write-Host "done"
